# clubPerson.pptx (diagrams/uml/aggregation) — slide 29 of the source deck.
#
# The canonical diff for this commit is dominated by a whole-file
# re-serialisation (notes master/notes slide dropped, slide-number /
# date field GUIDs refreshed, shape ids renumbered, endParaRPr lang
# flipped en-US -> en-SG, theme/layout boilerplate refreshed to a newer
# Office template, etc.) — none of that is an addressable, semantic
# edit reachable through the PowerPoint object model; it is what the
# authoring app regenerates on its own whenever it resaves a deck.
#
# The one genuine, intentional content edit baked into this commit is
# that the slide no longer carries the stray, empty "Title" placeholder
# shape (it had no text and was left over from an earlier layout
# change) — so we remove it here.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Title 19") {
        $sh.Delete()
    }
}
